$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values and clear trailing cells that are removed in the new (shorter) rows
$ws.Range("B2").Value = 0.1264008423207837
$ws.Range("C2").Value = 0.6808472755916881
$ws.Range("D2").Value = 0.04988061626763002
$ws.Range("E2").Value = 0.5970894115568507
$ws.Range("F2").Value = 0.3895648707313746
$ws.Range("G2").Value = 0.3749198787210216
$ws.Range("H2").Value = 0.4613573173527261
$ws.Range("I2").Value = 0.6892627280777406
$ws.Range("J2").Value = 0.1844338218533179
$ws.Range("K2").Value = 0.4715052544735016

$ws.Range("B3").Value = 0.7032752552246967
$ws.Range("C3").Value = 0.0524841558300787
$ws.Range("D3").Value = 0.5750217259028355
$ws.Range("E3").Value = 0.3835828904270196
$ws.Range("F3").Value = 0.3695766285386105
$ws.Range("G3").Value = 0.4509887839823598
$ws.Range("H3").Value = 0.6807854059541167
$ws.Range("I3").Value = 0.1766615135465071
$ws.Range("J3").Value = 0.4628812809405329
$ws.Range("K3").Value = 0.2074249537672726

$ws.Range("B4").Value = 0.0739087272872988
$ws.Range("C4").Value = 0.6951995747020479
$ws.Range("D4").Value = 0.2943885460132365
$ws.Range("E4").Value = 0.3398496246900327
$ws.Range("F4").Value = 0.4596149040122699
$ws.Range("G4").Value = 0.6548452325286815
$ws.Range("H4").Value = 0.1538829505182796
$ws.Range("I4").Value = 0.4492632457919151
$ws.Range("J4").Value = 0.1889041565820968
$ws.Range("K4").Value = 0.5706702220727796

$ws.Range("B5").Value = 0.6528789423816584
$ws.Range("C5").Value = 0.2651840721575033
$ws.Range("D5").Value = 0.3368290248851115
$ws.Range("E5").Value = 0.4407536204007895
$ws.Range("F5").Value = 0.6345141014634773
$ws.Range("G5").Value = 0.1393527950840318
$ws.Range("H5").Value = 0.4326222002996472
$ws.Range("I5").Value = 0.1713203111533466
$ws.Range("J5").Value = 0.5541963385427369
$ws.Range("K5").Value = 0.1437698493309027

$ws.Range("B6").Value = 0.6053818127754134
$ws.Range("C6").Value = 0.4122000866690486
$ws.Range("D6").Value = 0.2496603340877904
$ws.Range("E6").Value = 0.6592080140502106
$ws.Range("F6").Value = 0.1487321986403278
$ws.Range("G6").Value = 0.3778114016882561
$ws.Range("H6").Value = 0.1524262202646768
$ws.Range("I6").Value = 0.5393323377276911
$ws.Range("J6").Value = 0.115058138701532
$ws.Range("K6").Value = 0.4067718394308724

$ws.Range("B7").Value = 0.8628949586592991
$ws.Range("C7").Value = 0.2967710363001488
$ws.Range("D7").Value = 0.4189247832594023
$ws.Range("E7").Value = 0.1846772797061906
$ws.Range("F7").Value = 0.3757606442486632
$ws.Range("G7").Value = 0.07585798082864662
$ws.Range("H7").Value = 0.5119329433524077
$ws.Range("I7").Value = 0.08876908850380663
$ws.Range("J7").Value = 0.3633745487175398
$ws.Range("K7").ClearContents()

$ws.Range("B8").Value = 0.6090966232236873
$ws.Range("C8").Value = 0.5522135229949265
$ws.Range("D8").Value = 0.005598857889999004
$ws.Range("E8").Value = 0.4039548830192304
$ws.Range("F8").Value = 0.1118832920210401
$ws.Range("G8").Value = 0.4743913731481941
$ws.Range("H8").Value = 0.08322674941644539
$ws.Range("I8").Value = 0.3675498776562884
$ws.Range("J8").ClearContents()

$ws.Range("B9").Value = 0.7878040141027678
$ws.Range("C9").Value = 0.09027759876430858
$ws.Range("D9").Value = 0.2583545163855133
$ws.Range("E9").Value = 0.1218370348802827
$ws.Range("F9").Value = 0.489756542847739
$ws.Range("G9").Value = 0.04437841445902233
$ws.Range("H9").Value = 0.3517040686291025
$ws.Range("I9").ClearContents()

$ws.Range("B10").Value = 0.4013017852456914
$ws.Range("C10").Value = 0.3754432907967085
$ws.Range("D10").Value = -0.04062710656928412
$ws.Range("E10").Value = 0.5187154933129405
$ws.Range("F10").Value = 0.08012128691392592
$ws.Range("G10").Value = 0.3203764222454754
$ws.Range("H10").ClearContents()

$ws.Range("B11").Value = 0.6222684682008229
$ws.Range("C11").Value = -0.02297123903139461
$ws.Range("D11").Value = 0.4240932542019461
$ws.Range("E11").Value = 0.112338675162406
$ws.Range("F11").Value = 0.3327645480731927
$ws.Range("G11").ClearContents()

$ws.Range("B12").Value = 0.2167051203848173
$ws.Range("C12").Value = 0.5091174976711597
$ws.Range("D12").Value = -0.004145903195608092
$ws.Range("E12").Value = 0.3478698197250452
$ws.Range("F12").ClearContents()

$ws.Range("B13").Value = 0.6739775747052469
$ws.Range("C13").Value = 0.009391369052308113
$ws.Range("D13").Value = 0.2848969007350822
$ws.Range("E13").ClearContents()

$ws.Range("B14").Value = 0.2632404109177161
$ws.Range("C14").Value = 0.3842149509171186
$ws.Range("D14").ClearContents()

$ws.Range("B15").Value = 0.4282746421565676
$ws.Range("C15").ClearContents()

$ws.Range("B16").ClearContents()

